$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "District of Columbia" rows (South region block row 24,
# South Atlantic region block row 75). Delete the lower one first so the
# first row index is not shifted before we get to delete it.
$ws.Rows.Item(75).Delete()
$ws.Rows.Item(24).Delete()

$ws.Rows.Item(74).Select()
